# Edit: rename headers on existing sheets and add a new "PO Forecast" sheet
# with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# --- Update header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$header = $wsForecast.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Data rows
$data = @(
    @(45144.99999999999, 32, 18.91146929087808, 44.66114982910397),
    @(45151.99999999999, 30, 17.77875076595893, 43.00219181350299),
    @(45165.99999999999, 27, 14.89716820603103, 39.95668319832985),
    @(45186.99999999999, 22, 9.822226400844846, 35.4732350512764),
    @(45193.99999999999, 20, 8.372233051861999, 32.97346911429479),
    @(45200.99999999999, 18, 5.05837194892771, 30.41326640416146),
    @(45207.99999999999, 17, 4.014314127900851, 28.87196912574225),
    @(45214.99999999999, 15, 2.031326436264083, 28.61167432718508),
    @(45221.99999999999, 13, 1.608926898754556, 25.97255274312321),
    @(45228.99999999999, 12, -0.4405850068052672, 24.12494921424328),
    @(45235.99999999999, 10, -2.894658880772301, 21.72999803468937),
    @(45242.99999999999, 8, -3.127930755001636, 21.18986481932131)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
